# Append a newly-scraped Lancers listing (captured 2025-09-30 06:35:16) to the
# top of the "ランサーズ" data table, pushing the previously-seen rows down by
# one, refresh every row's "取得日時" timestamp to the new scrape time, widen
# column B to fit the new (longer) title, and rewire the per-row hyperlinks so
# they keep pointing at the correct "URL" cell after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "ランサーズ" sheet

# --- 1. Make room for the new entry -----------------------------------------
# Inserting a whole row at position 2 shifts the 5 existing data rows
# (formerly rows 2-6) down to rows 3-7 and bumps the sheet dimension to
# A1:H7 automatically.
$ws.Rows.Item(2).Insert()

# --- 2. Column B needs to be a bit wider for the new, longer title ----------
# Excel's ColumnWidth property measures in "characters" of the workbook's
# normal font, which gets padded by ~5/6 of a character when the stored
# OOXML <col width="..."> is computed. Subtracting that padding here makes
# the persisted width come out to exactly 52.
$ws.Columns.Item(2).ColumnWidth = 52 - (5 / 6)

# --- 3. Write the new row of data -------------------------------------------
$ws.Range("A2").Value = "2025-09-30 06:35:16"
$ws.Range("B2").Value = "【限定タスク】SIM AI の Google 認証ログイン機能の「最終調整」のみ代行(環境構築済み)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5403583"
$ws.Range("G2").Value = 298
$ws.Range("H2").Value = "🔥AI,Ai"

# --- 4. Refresh the scrape timestamp on every pre-existing row too ---------
$ws.Range("A3").Value = "2025-09-30 06:35:16"
$ws.Range("A4").Value = "2025-09-30 06:35:16"
$ws.Range("A5").Value = "2025-09-30 06:35:16"
$ws.Range("A6").Value = "2025-09-30 06:35:16"
$ws.Range("A7").Value = "2025-09-30 06:35:16"

# --- 5. Rebuild the hyperlinks on column F ----------------------------------
# Row-insert does not relocate the existing hyperlink anchors, so drop the
# whole collection and re-add one hyperlink per data row (F2:F7) in order;
# this also regenerates the relationship ids (rId1..rId6) in the same
# F2->rId1, F3->rId2, ... order used by the workbook.
$ws.Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5403583",
    "https://www.lancers.jp/work/detail/5403166",
    "https://www.lancers.jp/work/detail/5403527",
    "https://www.lancers.jp/work/detail/5403384",
    "https://www.lancers.jp/work/detail/5403072",
    "https://www.lancers.jp/work/detail/5399347"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($cell, $urls[$i], [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value) | Out-Null
    # Re-apply the named "Hyperlink" cell style so the cell keeps using the
    # workbook's existing style record instead of a newly synthesized one.
    $cell.Style = "Hyperlink"
}
